# Add team Wins/Losses/Ties columns (AD:AF) to the OAK_1994 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the existing header formatting (bold, centered, thin border - style
# used by A1:AC1) by copying it onto the three new header cells, then set
# their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Every player row (2-45) gets the team's 1994 record: 51 wins, 63 losses,
# 0 ties.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 51  # AD
    $ws.Cells.Item($r, 31).Value2 = 63  # AE
    $ws.Cells.Item($r, 32).Value2 = 0   # AF
}
